$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (bold, centered, thin-bordered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New data columns I (I0) and J (IF) for rows 2-9.
$iValues = @(7, 8, 5, 8, 6, 5, 9, 7)
$jValues = @(8, 8, 7, 8, 6, 5, 9, 7)

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r - 2]
    $ws.Cells.Item($r, 10).Value = $jValues[$r - 2]
}
